# Auto-generated data-driven update for 'Recommandations' and 'Top_YTD' sheets
# matching the GitHub Actions automated BRVM data refresh.

$sheet1Rows = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3410.86, 115.08, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2885, 770, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2820, 700, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2419.64, 605.27, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 2388.92, 606.72, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1434.37, 359.82, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1433.94, 346.24, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 904.28, 229.76, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 575.14, 143.66, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 574.53, 143.59, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 564.65, 141.12, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 504.41, 123.7, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 430.47, 108.19, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 2, 423.48, 213.74, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 384.76, 95.46, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                (**)", 0, 1, 235.98, 235.98, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 1, 218.47, 218.47, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE         (**)", 0, 1, 202.47, 202.47, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 193.64, 193.64, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 3, 0, 21.96, 7.08, "🟢 Achat", "✅ Renforcer"),
    @("NEI-CEDA CI (NEIC)", 3, 0, 16.94, 7.14, "🟢 Achat", "✅ Renforcer"),
    @("SOLIBRA CI (SLBC)", 2, 0, 14.98, 7.49, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 2, 0, 14.48, 6.98, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 1, 0, 7.5, 7.5, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 2, 1, 7.32, -7.5, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 1, 0, 7.25, 7.25, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 1, 0, 5.68, 5.68, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 1, 0, 5.6, 5.6, "🟡 Observer", "➖ Neutre"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 1, 0, 2.76, 2.76, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 1, 0, 2.59, 2.59, "🟡 Observer", "➖ Neutre"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 1, 1, -0.01, -7.5, "🟡 Observer", "👀 À surveiller"),
    @("TOTALENERGIES MARKETING SN (TTLS)", 0, 1, -0.99, -0.99, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -1.67, -1.67, "🟡 Observer", "➖ Neutre"),
    @("SAPH CI (SPHC)", 1, 1, -2, 2.53, "🟡 Observer", "👀 À surveiller"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -2.03, -2.03, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA BF (BOABF)", 0, 1, -3.31, -3.31, "🟡 Observer", "➖ Neutre"),
    @("SOGB CI (SOGC)", 0, 1, -3.83, -3.83, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, -4.55, -4.55, "🟡 Observer", "➖ Neutre"),
    @("PALM CI (PALC)", 0, 1, -7.45, -7.45, "🟡 Observer", "➖ Neutre"),
    @("NESTLE CI (NTLC)", 0, 1, -7.48, -7.48, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA ML (BOAM)", 0, 2, -8.84, -7.47, "🟡 Observer", "➖ Neutre"),
    @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 3, -9.68, -4.63, "🔴 Vente", "⚠️ Risque de décrochage"),
    @("FILTISAC CI (FTSC)", 0, 4, -18.96, -7.5, "🔴 Vente", "⚠️ Risque de décrochage"),
)

$sheet2Rows = @(
    @("BRVM - SERVICES PUBLICS", 10311390.25),
    @("NEI-CEDA CI", 453259.28),
    @("AIR LIQUIDE CI", 419804),
    @("BRVM - AUTRES SECTEURS", 246797),
    @("BRVM - DISTRIBUTION", 235574.04),
    @("BRVM - TRANSPORT", 44122.16),
    @("BRVM - AGRICULTURE", 44030.68),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 11170.56),
    @("BRVM-PRESTIGE", 3432.01),
    @("BRVM - FINANCES", 3423.22),
)

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

$r = 2
foreach ($row in $sheet1Rows) {
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

$r = 2
foreach ($row in $sheet2Rows) {
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
